$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Updated values in column C (rows 2-44) ---
$ws.Range("C2").Value = 2934.187009790061
$ws.Range("C3").Value = 2870.311589353206
$ws.Range("C4").Value = 1460.056109840828
$ws.Range("C5").Value = 6128.19547247793
$ws.Range("C6").Value = 4547.50930098406
$ws.Range("C8").Value = 1280.225469721551
$ws.Range("C9").Value = 2983.242707849043
$ws.Range("C10").Value = 2898.942214704482
$ws.Range("C11").Value = 1503.870423231357
$ws.Range("C12").Value = 1955.461557360978
$ws.Range("C13").Value = 6336.709213679884
$ws.Range("C15").Value = 1263.452411343738
$ws.Range("C16").Value = 3083.80337578809
$ws.Range("C17").Value = 2965.153206179127
$ws.Range("C18").Value = 1577.487171555845
$ws.Range("C19").Value = 6711.616186806423
$ws.Range("C21").Value = 1291.622214254295
$ws.Range("C23").Value = 6911.59200404802
$ws.Range("C24").Value = 1657.651524528445
$ws.Range("C25").Value = 3156.723844635973
$ws.Range("C26").Value = 1291.415042301529
$ws.Range("C28").Value = 7200.731056811853
$ws.Range("C29").Value = 1716.389195271215
$ws.Range("C30").Value = 3212.740625904757
$ws.Range("C31").Value = 7449.08671983612
$ws.Range("C32").Value = 3252.634165082374
$ws.Range("C33").Value = 1775.027517189621
$ws.Range("C34").Value = 1338.716747746975
$ws.Range("C36").Value = 7580.275568826287
$ws.Range("C37").Value = 3314.741082534716
$ws.Range("C38").Value = 1836.014008604312
$ws.Range("C39").Value = 1384.519227335143
$ws.Range("C41").Value = 7633.969039669125
$ws.Range("C42").Value = 3382.563653843273
$ws.Range("C43").Value = 3242.636921959078
$ws.Range("C44").Value = 1431.756130822538

